# Data Dictionary workbook update: add CER_Probable sheet, rename existing
# sheets with the CER_ prefix, and tidy up the Committed table's formatting.

$wb = $excel.ActiveWorkbook

$wsApproved = $wb.Worksheets.Item("Approved")
$wsCommitted = $wb.Worksheets.Item("Committed")

$wsApproved.Name = "CER_Approved"
$wsCommitted.Name = "CER_Committed"

# --- Clean up the stray duplicate "fill" styling that used to mark rows 7-8
#     of CER_Committed; make them match the rest of the table's styling.
$wsCommitted.Range("A2:D2").Copy()
$wsCommitted.Range("A7:D7").PasteSpecial(-4122)
$wsCommitted.Range("A8:D8").PasteSpecial(-4122)

$wsCommitted.Range("E6").Copy()
$wsCommitted.Range("E7:E8").PasteSpecial(-4122)
$wsCommitted.Range("E7:E8").VerticalAlignment = -4107

# --- Add the new CER_Probable sheet after CER_Committed
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsProbable = $wb.Worksheets.Add($null, $lastSheet)
$wsProbable.Name = "CER_Probable"

# Copy formats (cell styles + column widths) from CER_Committed so the new
# sheet matches the look of its siblings.
$wsCommitted.Range("A1:E6").Copy()
$wsProbable.Range("A1").PasteSpecial(-4122)
foreach ($col in @("A", "B", "C", "D", "E")) {
    $wsProbable.Columns($col).ColumnWidth = $wsCommitted.Columns($col).ColumnWidth
}

# --- Header row
$wsProbable.Range("A1").Value = "No"
$wsProbable.Range("B1").Value = "Column"
$wsProbable.Range("C1").Value = "Original Column"
$wsProbable.Range("D1").Value = "Type"
$wsProbable.Range("E1").Value = "Description"

# --- Data rows
$wsProbable.Range("A2").Value = 1
$wsProbable.Range("B2").Value = "project_name"
$wsProbable.Range("C2").Value = "Project Name"
$wsProbable.Range("D2").Value = "String"
$wsProbable.Range("E2").Value = "Name of the large-scale renewable project/site (may include stage)."

$wsProbable.Range("A3").Value = 2
$wsProbable.Range("B3").Value = "state"
$wsProbable.Range("C3").Value = "State"
$wsProbable.Range("D3").Value = "String"
$wsProbable.Range("E3").Value = "Australian state/territory where the project is located."

$wsProbable.Range("A4").Value = 3
$wsProbable.Range("B4").Value = "capacity_mw"
$wsProbable.Range("C4").Value = "MW Capacity"
$wsProbable.Range("D4").Value = "Numeric"
$wsProbable.Range("E4").Value = "Nameplate/announced capacity for the project."

$wsProbable.Range("A5").Value = 4
$wsProbable.Range("B5").Value = "fuel_source"
$wsProbable.Range("C5").Value = "Fuel Source"
$wsProbable.Range("D5").Value = "String"
$wsProbable.Range("E5").Value = "Primary energy resource/technology. Typical values: Solar, Wind, Bioenergy."

$wsProbable.Range("A6").Value = 5
$wsProbable.Range("B6").Value = "project_stage"
$wsProbable.Range("C6").Value = "(derived)"
$wsProbable.Range("D6").Value = "String"
$wsProbable.Range("E6").Value = "Status of the project (probable)"

$wsProbable.Range("E9").Select() | Out-Null

# --- Make CER_Approved the active tab again (it was the first sheet)
$wsApproved.Activate()
$wsApproved.Range("D12").Select() | Out-Null
